$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in day 9 (row 10) data. Order matches how the shared strings table
# grows (image_url, title, person, text) so new entries land at the same
# indices as the target workbook.
$ws.Range("E10").Value = "https://cdn.humanresourcesmanager.de/app/uploads/2019/09/achtsamkeit-am-arbeitsplatz.jpg"
$ws.Range("B10").Value = "Auf meinen Körper hören"
$ws.Range("D10").Value = "kein falscher Stolz"
$ws.Range("C10").Value = "Ich hatte heute Nachmittag so stark Kopfschmerzen, dass ich mich auf nichts mehr konzentrieren konnte. Anstatt mich durch den Nachmittag zu mogeln und schlussendlich bei der Arbeit doch nichts wirklich anzugehen, habe ich mich am Nachmittag krank abgemeldet. Es kam falscher Stolz mit, weil ich schon mehr als ein Jahr lang keine Stunde mehr krank gefehlt habe. Habe es dann trotzdem gemacht und auf meinen Körper gehört."

# Update the active cell selection on the sheet
$ws.Range("D6").Select()
